{"js": "// Remove the \"Stefanos Georgiou, Stamatia Rizou, ...\" reference paragraph\n// (a submitted-for-review entry under Research Publications > Conference\n// Publications) that immediately precedes the \"Invited Talks\" section.\nconst body = context.document.body;\n\n// Locate the paragraph via a distinctive text fragment rather than a\n// hard-coded index, so the edit is resilient to unrelated changes\n// elsewhere in the document.\nconst results = body.search(\"Stamatia Rizou\", { matchCase: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const match = results.items[0];\n  const paragraphs = match.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"Stefanos Georgiou, Stamatia Rizou, ...\" reference paragraph\n# (a submitted-for-review entry under Research Publications > Conference\n# Publications) that immediately precedes the \"Invited Talks\" section.\n$d = $word.ActiveDocument\n\n# Locate the paragraph via a distinctive text fragment rather than a\n# hard-coded index, so the edit is resilient to unrelated changes\n# elsewhere in the document.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Stamatia Rizou\"\n$found = $find.Execute()\n\nif ($found) {\n    $para = $rng.Paragraphs(1)\n    $para.Range.Delete()\n}\n"}
